$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (rows 2-10), columns A-G
# A = Fixture, B = Pick, C = AI_Confidence, D = OLBG_Confidence,
# E = Oddspedia_Confidence, F = Odds, G = Result

$data = @(
    @("Borussia Mönchengladbach - Bayern Munich : 0:0'", "Bayern Munich", 73, $null, 98, 1.28, $null),
    @("Raja Club Athletic  - Olympique Dcheira: 18:00", "Raja Club Athletic", 70, 87, 89, 1.67, $null),
    @("Chelsea FC  - Sunderland AFC: 1:1'", "Chelsea FC", 67, $null, 77, 1.45, $null),
    @("Inter Club d'Escaldes  - FC Ordino: -:-'", "Inter Club d'Escaldes", 64, 59, $null, 2.5, $null),
    @("SL Benfica  - FC Arouca: 20:30", "SL Benfica", 61, 100, $null, 1.18, $null),
    @("Levski Sofia  - Dobrudzha Dobrich: -:-'", "Levski Sofia", 59, 78, $null, 1.18, $null),
    @("Borussia Dortmund  - 1.FC Köln: -:-'", "Borussia Dortmund", 57, 96, 100, 1.42, $null),
    @("Inter Miami CF ✓ - Nashville SC: 3:1", "Inter Miami CF", 56, 59, $null, 2.5, "✓"),
    @("Struga Trim & Lum  - AP Brera Strumica: -:-'", "Struga Trim & Lum", 53, $null, 94, 1.83, $null)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    if ($row[3] -eq $null) {
        $ws.Cells.Item($rowIndex, 4).Value = ""
    } else {
        $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    }
    if ($row[4] -eq $null) {
        $ws.Cells.Item($rowIndex, 5).Value = ""
    } else {
        $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    }
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    if ($row[6] -eq $null) {
        $ws.Cells.Item($rowIndex, 7).Value = ""
    } else {
        $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    }
    $rowIndex++
}
